$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range('G9').Value = 'Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range('G18').Value = 'eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G19').Value = 'neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G20').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G21').Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range('G24').Value = 'NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G28').Value = 'Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Range('G37').Value = 'eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G38').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G39').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G40').Value = 'hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range('G43').Value = 'NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G47').Value = 'maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('G56').Value = 'eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G57').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G58').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G59').Value = 'Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G66').Value = 'Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range('G75').Value = 'eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G76').Value = 'neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G77').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G78').Value = 'Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G81').Value = 'NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G85').Value = 'Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range('G94').Value = 'eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G95').Value = 'neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G96').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
$ws.Range('G97').Value = 'Mohammedeltanany@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G100').Value = 'NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Range('G104').Value = 'maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range('G113').Value = 'eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G115').Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg'
